$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current header row, shifting
# everything (including the two hyperlinked cells) down by one row.
$ws.Rows("1:1").Insert()

# Populate the new row 1 with the mandatory reference columns.
$ws.Range("A1").Value = "Reference"
$ws.Range("B1").Value = "UNAM"

# The row insert does not automatically re-anchor existing hyperlinks in
# this runtime, so rebuild them at their new (shifted-down) locations,
# keeping the existing cell text ("Michoacán" / "Veracruz") as-is.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B18"), "https://en.wikipedia.org/wiki/Michoac%C3%A1n", "", "Michoacán")
$ws.Hyperlinks.Add($ws.Range("B20"), "https://en.wikipedia.org/wiki/Veracruz", "", "Veracruz")

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink"
# style; restore the original (non-hyperlink) formatting on those cells.
$ws.Range("A18").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("A20").Copy()
$ws.Range("B20").PasteSpecial(-4122)

# Drop the now-unused "Hyperlink" cell style that Excel auto-registered.
$wb.Styles.Item("Hyperlink").Delete()

# Match the new selection state.
$ws.Range("B1").Select()
